# Apply the evaluation updates to the QuantitativeMetrics sheet of UC1_TC2.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Row 6 - "Runtime without error": value changes from yes -> no, and a note is added
$ws.Range("B6").Value = "no"
$ws.Range("C6").Value = "Missing initial redirect"

# Row 7 - "Assertion validity": value + note are cleared out
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Row 12 - "Code BLEU": updated score and updated breakdown note
$ws.Range("B12").Value = 0.1736050408391118
$ws.Range("C12").Value = "{'codebleu': 0.17360504083911182, 'ngram_match_score': 0.038139835860325454, 'weighted_ngram_match_score': 0.08106510110122213, 'syntax_match_score': 0.426497277676951, 'dataflow_match_score': 0.14871794871794872}"

# Update the active selection on the sheet to C7 (matches the new cursor position)
$ws.Range("C7").Select()

$wb.Save()
